# Add data for 2022-06-08
# - Rename sheet to reflect new "through" date
# - Update the "2022 (through ...)" header label shared string
# - Update June (row 6) and Total (row 14) values in the "2022" column (I)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Through 2022-05-31"

# Update the column header text for the 2022 series
$ws.Range("I1").Value = "2022 (through 05-31)"

# Update the June value for 2022
$ws.Range("I6").Value = 114

# Update the Total value for 2022
$ws.Range("I14").Value = 665
